# Fill in the first empty row of the time-tracking table with a new entry.
$d = $word.ActiveDocument

$table = $d.Tables.Item(1)
$row = $table.Rows.Item(13)

$row.Cells.Item(1).Range.Text = "10.10.2022"
$row.Cells.Item(2).Range.Text = "0,25"
$row.Cells.Item(3).Range.Text = "Palaveri"
